{"js": "// Apply the tracked text replacements (date header + multiplication problems).\nconst body = context.document.body;\nconst replacements = [\n  [\"2024-05-10 Friday\", \"2024-05-11 Saturday\"],\n  [\"84\u00d737=3108\", \"51\u00d740=2040\"],\n  [\"61\u00d774=4514\", \"48\u00d790=4320\"],\n  [\"49\u00d726=1274\", \"17\u00d734=578\"],\n  [\"71\u00d788=6248\", \"81\u00d766=5346\"],\n  [\"54\u00d726=1404\", \"97\u00d742=4074\"],\n  [\"97\u00d723=2231\", \"97\u00d793=9021\"],\n  [\"18\u00d777=1386\", \"85\u00d772=6120\"],\n  [\"82\u00d799=8118\", \"90\u00d751=4590\"],\n  [\"22\u00d768=1496\", \"58\u00d741=2378\"],\n  [\"93\u00d730=2790\", \"39\u00d770=2730\"],\n  [\"99\u00d766=6534\", \"86\u00d791=7826\"],\n  [\"12\u00d785=1020\", \"74\u00d740=2960\"],\n  [\"77\u00d778=6006\", \"38\u00d757=2166\"],\n  [\"27\u00d759=1593\", \"18\u00d742=756\"],\n  [\"64\u00d716=1024\", \"54\u00d757=3078\"],\n  [\"73\u00d775=5475\", \"47\u00d778=3666\"],\n  [\"22\u00d745=990\", \"21\u00d773=1533\"],\n  [\"38\u00d777=2926\", \"53\u00d739=2067\"],\n  [\"76\u00d735=2660\", \"54\u00d794=5076\"],\n  [\"74\u00d747=3478\", \"62\u00d740=2480\"],\n  [\"21\u00d734=714\", \"41\u00d771=2911\"],\n  [\"37\u00d755=2035\", \"42\u00d750=2100\"],\n  [\"31\u00d769=2139\", \"44\u00d725=1100\"],\n  [\"79\u00d793=7347\", \"43\u00d763=2709\"],\n  [\"26\u00d719=494\", \"86\u00d720=1720\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the tracked text replacements (date header + multiplication problems).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@('2024-05-10 Friday', '2024-05-11 Saturday')\n    ,@('84\u00d737=3108', '51\u00d740=2040')\n    ,@('61\u00d774=4514', '48\u00d790=4320')\n    ,@('49\u00d726=1274', '17\u00d734=578')\n    ,@('71\u00d788=6248', '81\u00d766=5346')\n    ,@('54\u00d726=1404', '97\u00d742=4074')\n    ,@('97\u00d723=2231', '97\u00d793=9021')\n    ,@('18\u00d777=1386', '85\u00d772=6120')\n    ,@('82\u00d799=8118', '90\u00d751=4590')\n    ,@('22\u00d768=1496', '58\u00d741=2378')\n    ,@('93\u00d730=2790', '39\u00d770=2730')\n    ,@('99\u00d766=6534', '86\u00d791=7826')\n    ,@('12\u00d785=1020', '74\u00d740=2960')\n    ,@('77\u00d778=6006', '38\u00d757=2166')\n    ,@('27\u00d759=1593', '18\u00d742=756')\n    ,@('64\u00d716=1024', '54\u00d757=3078')\n    ,@('73\u00d775=5475', '47\u00d778=3666')\n    ,@('22\u00d745=990', '21\u00d773=1533')\n    ,@('38\u00d777=2926', '53\u00d739=2067')\n    ,@('76\u00d735=2660', '54\u00d794=5076')\n    ,@('74\u00d747=3478', '62\u00d740=2480')\n    ,@('21\u00d734=714', '41\u00d771=2911')\n    ,@('37\u00d755=2035', '42\u00d750=2100')\n    ,@('31\u00d769=2139', '44\u00d725=1100')\n    ,@('79\u00d793=7347', '43\u00d763=2709')\n    ,@('26\u00d719=494', '86\u00d720=1720')\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null,$true,$false,$false,$false,$false,$true,0,$false,$null,2) | Out-Null\n}\n"}
